$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 134 (shifts old rows 134..164 down to 135..165),
# copying formatting from the row that is currently row 134 so the date column
# (D) keeps its date number format style.
$ws.Rows.Item(134).Insert()

# Fill the newly inserted row 134 with the same "shape" of data as every other
# record in this table (same market/region/category/variety/quality/unit/
# origin/measure), only the date, volume and the three price columns differ.
$ws.Cells.Item(134, 1).Value = 7
$ws.Cells.Item(134, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(134, 3).Value = "Ñuble"
$ws.Cells.Item(134, 4).Value = 44543
$ws.Cells.Item(134, 5).Value = 16
$ws.Cells.Item(134, 6).Value = 100112003
$ws.Cells.Item(134, 7).Value = "Ajo"
$ws.Cells.Item(134, 8).Value = "Chino"
$ws.Cells.Item(134, 9).Value = "Primera"
$ws.Cells.Item(134, 10).Value = 60
$ws.Cells.Item(134, 11).Value = 18000
$ws.Cells.Item(134, 12).Value = 19000
$ws.Cells.Item(134, 13).Value = 18500
$ws.Cells.Item(134, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(134, 15).Value = "China"
$ws.Cells.Item(134, 16).Value = 1850
$ws.Cells.Item(134, 17).Value = 10
$ws.Cells.Item(134, 18).Value = "Hortaliza"
